$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data pattern to append: same sequence as rows 2-11 (1,2,3,4,5 repeated twice)
$values = 1,2,3,4,5,1,2,3,4,5

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 12 + $i
    $val = $values[$i]
    $ws.Cells.Item($row, 1).Value = $val
    $ws.Cells.Item($row, 2).Value = $val
}
